$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4 (Nandan) - replace the old note with the updated breakdown figures
$ws.Range("H4").Value = "5204-3600=1604"
$ws.Range("I4").Value = "1604-387=1217"
$ws.Range("J4").Value = 7755
$ws.Range("K4").Value = -4820

# Row 9 label - individual cost breakdown
$ws.Range("H9").Value = "Individual Cost"
$ws.Range("H9").Font.Bold = $true

# --- "Total" / "Balance" summary header for the per-person spend table ---
$ws.Range("J1").Value = "Total"
$ws.Range("J1").Font.Bold = $true
$ws.Range("K1").Value = "Balance"
$ws.Range("K1").Font.Bold = $true

# Row 9 value - individual cost formula note
$ws.Range("I9").Value = "11742/4=2935.5"

# Row 2 (Krushik)
$ws.Range("H2").Value = 387
$ws.Range("J2").Value = 1587
$ws.Range("K2").Value = 1348

# Row 3 (Likith)
$ws.Range("J3").Value = 1200
$ws.Range("K3").Value = 1735

# Row 5 (Nikhil) - fill in amount spent + totals
$ws.Range("G5").Value = 1200
$ws.Range("J5").Value = 1200
$ws.Range("K5").Value = 1735

# Row 6 - grand total of all amounts spent
$ws.Range("I6").Value = "Total"
$ws.Range("I6").Font.Bold = $true
$ws.Range("J6").Value = 11742

# New column I now has content - size it like Excel's autofit would
$ws.Columns.Item(9).ColumnWidth = 13.4

# Match the author's final selection
[void]$ws.Range("K10").Select()
